$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.627.56'
$ws.Range('E2').Value = '  +3.74%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.38'
$ws.Range('E3').Value = '  +2.08%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.00'
$ws.Range('E5').Value = '  +1.58%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.697'
$ws.Range('E6').Value = '  +1.92%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.24'
$ws.Range('E8').Value = '  +2.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.88'
$ws.Range('E9').Value = '  +10.29%  '

$ws.Range('E10').Value = '  +3.30%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0765'
$ws.Range('E11').Value = '  +3.49%  '

$ws.Range('E12').Value = '  +2.60%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.54'
$ws.Range('E13').Value = '  +7.89%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.799'
$ws.Range('E14').Value = '  +3.78%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.196.79'
$ws.Range('E15').Value = '  +2.00%  '

$ws.Range('E16').Value = '  +4.60%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.925.71'
$ws.Range('E17').Value = '  +3.00%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.560.89'
$ws.Range('E18').Value = '  +3.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.27'
$ws.Range('E19').Value = '  +2.06%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  +5.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '251.67'
$ws.Range('E21').Value = '  +3.29%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.26'
$ws.Range('E22').Value = '  +3.83%  '

$ws.Range('E23').Value = '  +5.05%  '

$ws.Range('E24').Value = '  +1.95%  '

$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('E26').Value = '  +1.44%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.80'
$ws.Range('E27').Value = '  +1.47%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.82'
$ws.Range('E28').Value = '  +3.36%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.77'
$ws.Range('E29').Value = '  +2.80%  '

$ws.Range('E30').Value = '  +1.90%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.56'
$ws.Range('E31').Value = '  +7.10%  '

$ws.Range('E32').Value = '  +3.81%  '

$ws.Range('E33').Value = '  +1.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.36'
$ws.Range('E34').Value = '  +5.58%  '

$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  -13.13%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0843'
$ws.Range('E37').Value = '  +14.91%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.875'
$ws.Range('E38').Value = '  +4.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.86'
$ws.Range('E39').Value = '  +47.26%  '

$ws.Range('E40').Value = '  +4.52%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.65'
$ws.Range('E41').Value = '  +10.98%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0229'
$ws.Range('E42').Value = '  +5.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.19'
$ws.Range('E43').Value = '  -1.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.11'
$ws.Range('E44').Value = '  +3.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.339.50'
$ws.Range('E45').Value = '  +2.74%  '

$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.54'
$ws.Range('E46').Value = '  +6.76%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  +1.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0815'
$ws.Range('E48').Value = '  +2.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.79'
$ws.Range('E49').Value = '  +2.58%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.44'
$ws.Range('E50').Value = '  +3.53%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.07'
$ws.Range('E51').Value = '  +2.64%  '
